$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 3 first: since row 2 (above it) already carries the
# plain data-row formatting (only column D is date-formatted), the new blank row
# inherits that same formatting instead of the bold header style that a plain
# insert at row 2 would copy down from row 1.
$ws.Rows.Item(3).Insert()

# Move the old row 2's data down into the now-blank, correctly-formatted row 3 -
# this is the shift the diff performs on every existing record (old row N -> row N+1).
$ws.Range("A2:R2").Cut($ws.Range("A3"))

# Populate the now-empty row 2 with the new weekly record.
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 'Vega Monumental Concepción'
$ws.Range("C2").Value = 'Bíobío'
$ws.Range("D2").Value = 44496
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112012
$ws.Range("G2").Value = 'Espinaca'
$ws.Range("H2").Value = 'Sin especificar'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 650
$ws.Range("L2").Value = 700
$ws.Range("M2").Value = 671
$ws.Range("N2").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O2").Value = 'Región de Ñuble'
$ws.Range("P2").Value = 671
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 'Hortaliza'
